$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1080.75
$ws.Range("I12").Value = 329.2
$ws.Range("K12").Value = 329.2
$ws.Range("M12").Value = -159.2
$ws.Range("H43").Value = 9593
$ws.Range("J43").Value = 10166.777
$ws.Range("L43").Value = 10166.777
$ws.Range("N43").Value = -10304.777
$ws.Range("H81").Value = 136000
$ws.Range("J81").Value = 136000
$ws.Range("L81").Value = 136000
$ws.Range("N81").Value = -137996
$ws.Range("H84").Value = 136000
$ws.Range("J84").Value = 136000
$ws.Range("L84").Value = 408000
$ws.Range("N84").Value = -417984
$ws.Range("H106").Value = 3704943.2
$ws.Range("I106").Value = 4167987.8
$ws.Range("K106").Value = 4167987.8
$ws.Range("M106").Value = -4167356.8
$ws.Range("H113").Value = 8000.8
$ws.Range("I113").Value = 6502.5
$ws.Range("J113").Value = 8999.666999999999
$ws.Range("K113").Value = 6502.5
$ws.Range("L113").Value = 8999.666999999999
$ws.Range("M113").Value = -3248.5
$ws.Range("N113").Value = -15507.667
$ws.Range("H116").Value = 21036.334
$ws.Range("I116").Value = 22790.875
$ws.Range("J116").Value = 7000
$ws.Range("K116").Value = 22790.875
$ws.Range("L116").Value = 7000
$ws.Range("M116").Value = -19348.875
$ws.Range("N116").Value = -13884
$ws.Range("H132").Value = 2393
$ws.Range("I132").Value = 2702.375
$ws.Range("K132").Value = 8107.125
$ws.Range("M132").Value = -5577.125
$ws.Range("H137").Value = 3294.6155
$ws.Range("I137").Value = 4490
$ws.Range("J137").Value = 1900
$ws.Range("K137").Value = 13470
$ws.Range("L137").Value = 5700
$ws.Range("M137").Value = -10920
$ws.Range("N137").Value = -10800
$ws.Range("H138").Value = 3130.9797
$ws.Range("J138").Value = 3999.3635
$ws.Range("L138").Value = 11998.0905
$ws.Range("N138").Value = -22278.0905
$ws.Range("H141").Value = 3811.1177
$ws.Range("I141").Value = 3368.0625
$ws.Range("J141").Value = 10900
$ws.Range("K141").Value = 10104.1875
$ws.Range("L141").Value = 32700
$ws.Range("M141").Value = -4924.1875
$ws.Range("N141").Value = -43060

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4006.9167
$ws.Range("I61").Value = 3927.7144
$ws.Range("K61").Value = 3927.7144
$ws.Range("M61").Value = -3715.7144
$ws.Range("H118").Value = 14975
$ws.Range("J118").Value = 14975
$ws.Range("L118").Value = 14975
$ws.Range("N118").Value = -18289
$ws.Range("H132").Value = 2762.926
$ws.Range("I132").Value = 2583.96
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 7751.88
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -5221.88
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 4006.9167
$ws.Range("I136").Value = 3927.7144
$ws.Range("K136").Value = 11783.1432
$ws.Range("M136").Value = -9233.143199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2564.5789
$ws.Range("I20").Value = 2745.75
$ws.Range("K20").Value = 2745.75
$ws.Range("M20").Value = -2498.75
$ws.Range("H86").Value = 2223.2354
$ws.Range("I86").Value = 1946.1818
$ws.Range("J86").Value = 2731.1667
$ws.Range("K86").Value = 1946.1818
$ws.Range("L86").Value = 2731.1667
$ws.Range("M86").Value = -823.1818000000001
$ws.Range("N86").Value = -4977.1667
$ws.Range("H89").Value = 2223.2354
$ws.Range("I89").Value = 1946.1818
$ws.Range("J89").Value = 2731.1667
$ws.Range("K89").Value = 9730.909
$ws.Range("L89").Value = 13655.8335
$ws.Range("M89").Value = -4114.909
$ws.Range("N89").Value = -24887.8335
$ws.Range("H94").Value = 381.2143
$ws.Range("I94").Value = 433
$ws.Range("J94").Value = 70.5
$ws.Range("K94").Value = 433
$ws.Range("L94").Value = 70.5
$ws.Range("M94").Value = 18
$ws.Range("N94").Value = -972.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2195.8918
$ws.Range("I31").Value = 2016
$ws.Range("J31").Value = 2230.7097
$ws.Range("K31").Value = 2016
$ws.Range("L31").Value = 2230.7097
$ws.Range("M31").Value = -1721
$ws.Range("N31").Value = -2820.7097
$ws.Range("H34").Value = 2195.8918
$ws.Range("I34").Value = 2016
$ws.Range("J34").Value = 2230.7097
$ws.Range("K34").Value = 2016
$ws.Range("L34").Value = 2230.7097
$ws.Range("M34").Value = -1814
$ws.Range("N34").Value = -2634.7097
$ws.Range("H43").Value = 40900
$ws.Range("J43").Value = 40900
$ws.Range("L43").Value = 40900
$ws.Range("N43").Value = -41268
$ws.Range("H58").Value = 3270.8518
$ws.Range("I58").Value = 1871.6
$ws.Range("K58").Value = 1871.6
$ws.Range("M58").Value = -1668.6
$ws.Range("H101").Value = 40900
$ws.Range("J101").Value = 40900
$ws.Range("L101").Value = 40900
$ws.Range("N101").Value = -47390
$ws.Range("H104").Value = 79945
$ws.Range("J104").Value = 79945
$ws.Range("L104").Value = 79945
$ws.Range("N104").Value = -85187
$ws.Range("H107").Value = 46398.773
$ws.Range("I107").Value = 55996.445
$ws.Range("K107").Value = 55996.445
$ws.Range("M107").Value = -54076.445
$ws.Range("H132").Value = 4931.25
$ws.Range("I132").Value = 4716.905
$ws.Range("J132").Value = 6431.6665
$ws.Range("K132").Value = 14150.715
$ws.Range("L132").Value = 19294.9995
$ws.Range("M132").Value = -11620.715
$ws.Range("N132").Value = -24354.9995
$ws.Range("H134").Value = 5325.3335
$ws.Range("I134").Value = 5325.3335
$ws.Range("K134").Value = 15976.0005
$ws.Range("M134").Value = -13441.0005
$ws.Range("H136").Value = 3270.8518
$ws.Range("I136").Value = 1871.6
$ws.Range("K136").Value = 5614.799999999999
$ws.Range("M136").Value = -3064.799999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1357.875
$ws.Range("I17").Value = 1357.875
$ws.Range("K17").Value = 4073.625
$ws.Range("M17").Value = -3904.625
$ws.Range("H34").Value = 859.61536
$ws.Range("J34").Value = 1465.4286
$ws.Range("L34").Value = 4396.2858
$ws.Range("N34").Value = -4564.2858
$ws.Range("H55").Value = 4418.1665
$ws.Range("J55").Value = 5556.6665
$ws.Range("L55").Value = 16669.9995
$ws.Range("N55").Value = -17023.9995
$ws.Range("H86").Value = 224.375
$ws.Range("I86").Value = 199.6
$ws.Range("J86").Value = 265.66666
$ws.Range("K86").Value = 598.8
$ws.Range("L86").Value = 796.9999799999999
$ws.Range("M86").Value = 587.2
$ws.Range("N86").Value = -3168.99998
$ws.Range("H89").Value = 224.375
$ws.Range("I89").Value = 199.6
$ws.Range("J89").Value = 265.66666
$ws.Range("K89").Value = 1796.4
$ws.Range("L89").Value = 2390.99994
$ws.Range("M89").Value = 4131.6
$ws.Range("N89").Value = -14246.99994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H122").Value = 3036.4119
$ws.Range("I122").Value = 2582.25
$ws.Range("K122").Value = 7746.75
$ws.Range("M122").Value = -5296.75
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 15000
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -20060

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3107.6924
$ws.Range("I40").Value = 2710.111
$ws.Range("K40").Value = 2710.111
$ws.Range("M40").Value = -2574.111
$ws.Range("H93").Value = 34483840
$ws.Range("I93").Value = 47619890
$ws.Range("J93").Value = 1712.5
$ws.Range("K93").Value = 47619890
$ws.Range("L93").Value = 1712.5
$ws.Range("M93").Value = -47618642
$ws.Range("N93").Value = -4208.5
$ws.Range("H105").Value = 44995
$ws.Range("J105").Value = 44995
$ws.Range("L105").Value = 44995
$ws.Range("N105").Value = -51983
$ws.Range("H119").Value = 89533.336
$ws.Range("J119").Value = 89533.336
$ws.Range("L119").Value = 89533.336
$ws.Range("N119").Value = -99209.336
$ws.Range("H122").Value = 10497.25
$ws.Range("I122").Value = 6497
$ws.Range("J122").Value = 14497.5
$ws.Range("K122").Value = 19491
$ws.Range("L122").Value = 43492.5
$ws.Range("M122").Value = -17041
$ws.Range("N122").Value = -48392.5
$ws.Range("H132").Value = 10853.1
$ws.Range("I132").Value = 10345.471
$ws.Range("J132").Value = 13729.667
$ws.Range("K132").Value = 31036.413
$ws.Range("L132").Value = 41189.001
$ws.Range("M132").Value = -28506.413
$ws.Range("N132").Value = -46249.001
$ws.Range("H138").Value = 173885.5
$ws.Range("J138").Value = 173885.5
$ws.Range("L138").Value = 173885.5
$ws.Range("N138").Value = -184165.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4847.5557
$ws.Range("I122").Value = 5200.6665
$ws.Range("J122").Value = 3611.6667
$ws.Range("K122").Value = 15601.9995
$ws.Range("L122").Value = 10835.0001
$ws.Range("M122").Value = -13151.9995
$ws.Range("N122").Value = -15735.0001
